$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A rows 2..97 hold labels q1..q96. Shift each label down by one
# (q1 -> q0, q2 -> q1, ..., q96 -> q95) ahead of a dimension reduction.
for ($row = 2; $row -le 97; $row++) {
    $ws.Cells.Item($row, 1).Value = "q" + ($row - 2)
}
